# docs/achievement.xlsx - FlixelRL-731 unlock management implementation
# Splits the generic "nightmare" enemy achievement type into an "enemy" type
# with an explicit enemy-id param (NIGHTMARE / FIRE_DEMON / WATER_DEMON / ...),
# and widens/narrows the type (C) and param0 (D) columns to fit the new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width adjustments (C: type, D: param0) -------------------------
# ColumnWidth is expressed in characters; the engine quantises to whole
# pixels (MDW=7) same as Excel, so these are the closest achievable inputs
# to the target stored widths of 5.5 and 12.125 characters respectively.
$ws.Columns.Item(3).ColumnWidth = 34/7
$ws.Columns.Item(4).ColumnWidth = 80/7

# --- Row 18..25: nightmare achievements now use type "enemy" + named id ----
$enemyIds = @{
    18 = "NIGHTMARE"
    19 = "FIRE_DEMON"
    20 = "WATER_DEMON"
    21 = "EARTH_DEMON"
    22 = "WIND_DEMON"
    23 = "POISON_DEMON"
    24 = "SHADOW_DEMON"
    25 = "ICE_DEMON"
}

foreach ($row in 18..25) {
    $ws.Range("C$row").Value = "enemy"
    $ws.Range("D$row").Value = $enemyIds[$row]
}
